# Adds consolidated prod details: four new job rows (NQM AMC/CHL SFTP feeds)
# to the Hudson Advisor Capital Market BCP Data tracker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: Hudson_Daily_NQM_AMC_SFTP_ETL ---------------------------------
$ws.Range("A14").Value = "Hudson_Daily_NQM_AMC_SFTP_ETL"
$ws.Range("B14").Value = "yes"
$ws.Range("C14").Value = "Daily_AMC_SFTP_DataLoad"
$ws.Range("D14").Value = "C:\SSIS\Daily_Feed\Daily_AMC_SFTP_DataLoad.dtsx"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "NA"
$ws.Range("H14").Value = "NA"
$ws.Range("I14").Value = "NA"
$ws.Range("J14").Value = "Doesn" + [char]8217 + "t have BCp, Excel to DB"

# --- Row 15: Hudson_Daily_NQM_CHL_SFTP_ETL (download leg) ------------------
$ws.Range("A15").Value = "Hudson_Daily_NQM_CHL_SFTP_ETL"
$ws.Range("B15").Value = "yes"
$ws.Range("C15").Value = "Daily_CHL_SFTP_Data_Download"
$ws.Range("D15").Value = "C:\SSIS\Daily_Feed\Daily_CHL_SFTP_Data_Download.dtsx"

# --- Row 16: Hudson_Daily_NQM_CHL_SFTP_UPLOAD_ETL (upload leg) -------------
$ws.Range("A16").Value = "Hudson_Daily_NQM_CHL_SFTP_UPLOAD_ETL"
$ws.Range("B16").Value = "yes"
$ws.Range("C16").Value = "Daily_CHL_SFTP_Data_Upload"
$ws.Range("D16").Value = "C:\SSIS\Daily_Feed\Daily_CHL_SFTP_Data_Upload.dtsx"
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "NA"
$ws.Range("G16").Value = "NA"
$ws.Range("H16").Value = "NA"
$ws.Range("I16").Value = "NA"
$ws.Range("J16").Value = "Doesn" + [char]8217 + "t have BCP, generates txt and uploades to SFTP"

# --- Row 17: Hudson_Daily_NQM_ETL (parent/summary job) ---------------------
$ws.Range("A17").Value = "Hudson_Daily_NQM_ETL"
$ws.Range("B17").Value = "yes"
$ws.Range("C17").Value = "Daily_CHL_SFTP_Data_Download"
$ws.Range("D17").Value = "C:\SSIS\Daily_Feed\Daily_CHL_SFTP_Data_Download.dtsx"

# Widen column A so the longer job names (e.g. the CHL upload ETL name)
# remain fully visible, mirroring the prior best-fit autosize behaviour.
$ws.Columns.Item(1).AutoFit()

# Move / record the active selection the workbook was left at after the edit.
$ws.Range("C24").Select()
